# EPBDS-10072: introduce new keyword for all business enumeration properties
# in the file name processor. The "Tests" sample sheet documents how the
# generated JUnit test name is composed from the business-method name, so
# the sample header cell is updated to show the full generated name:
#   "<TestClass><MethodName> <MethodName>Test"
# i.e. the existing "TestSayHello" caption gets a new, bold "SayHelloTest"
# suffix appended (separated by a space), matching the naming convention
# produced by the updated file-name processor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

$headerCell = $ws.Range("C3")

# Rebuild the caption text: "Test" + "SayHello" (grey, not bold) + " " +
# "SayHelloTest" (black, bold).
$headerCell.Value = "TestSayHello SayHelloTest"

# "TestSayHello" (chars 1-12) keeps the original muted grey styling.
$greyPart = $headerCell.Characters(1, 12)
$greyPart.Font.Name = "Calibri"
$greyPart.Font.Size = 11
$greyPart.Font.Bold = $false
$greyPart.Font.Color = 8421504   # RGB(128,128,128) - theme0 tinted grey

# The separating space stays plain black.
$spacePart = $headerCell.Characters(13, 1)
$spacePart.Font.Name = "Calibri"
$spacePart.Font.Size = 11
$spacePart.Font.Bold = $false
$spacePart.Font.Color = 0

# New bold "SayHelloTest" suffix highlighting the generated test name.
$boldPart = $headerCell.Characters(14, 12)
$boldPart.Font.Name = "Calibri"
$boldPart.Font.Size = 11
$boldPart.Font.Bold = $true
$boldPart.Font.Color = 0

# The cursor/selection in the saved view moved from E9 up to E8.
$ws.Range("E8").Select()
